# Adds new competitor rows (83-96) to the "Competitors" sheet and
# refreshes the sheet-view / active-tab state to match the final
# workbook state (Competitors becomes the active / last-viewed sheet).

$wb = $excel.ActiveWorkbook

# Helpers that join text with NBSP (U+00A0) the same way the existing
# sheet does for "age<NBSP>NN", "NN<NBSP>kg" and "<NBSP>Country<NBSP>"
# values. NOTE: plain "+" concatenation between a numeric-looking
# string and a char silently does NUMERIC addition in this host, so
# these use the "-f" format operator instead, which always does a
# plain string join.
function NB2($a, $b) {
    $nbsp = [char]0x00A0
    return "{0}{1}{2}" -f $a, $nbsp, $b
}
function NB3($a, $b, $c) {
    $nbsp = [char]0x00A0
    return "{0}{1}{2}{3}{4}" -f $a, $nbsp, $b, $nbsp, $c
}

$competitors = $wb.Worksheets.Item("Competitors")
$countries   = $wb.Worksheets.Item("Countries")

# ---------------------------------------------------------------
# 1. Existing row 82 gained an "Id" (column A) and a "Sports"
#    (column J) value it was missing before.
# ---------------------------------------------------------------
$competitors.Cells.Item(82, 1).Value = 81                    # A82 - Id
$competitors.Cells.Item(82, 10).Value = "Rowing"              # J82

# ---------------------------------------------------------------
# 2. New competitor rows 83-96.
# ---------------------------------------------------------------
$rows = @(
  @{ Row=83; Id=82; B="Marta Walczykiewicz";   C="1 August 1987";      D=(NB2 "age" "32"); E="Women"; F=(NB2 "64" "kg");  G="Kalisz";           H="KTW Kalisz";            HWrap=$true;  J="Canoe Sprint";         K=(NB3 "" "Poland" "") },
  @{ Row=84; Id=83; B="Ariana Orrego";         C="25 September 1998";  D=(NB2 "age" "21"); E="Women"; F=(NB2 "49" "kg");  G="Lima";             H="Excalibur Gymnastics";  I="Gustavo Moure";        J="Gymnastics Artistic";  K=(NB3 "" "Peru" "") },
  @{ Row=85; Id=84; B="Hernán Viera";          C="16 January 1993 ";   D=(NB2 "age" "27"); E="Man";                       F="104 kg";                                                               J="Weightlifting";        K=(NB3 "" "Peru" "") },
  @{ Row=86; Id=85; B="Emilie Hegh Arntzen";   C="1 January 1994 ";    D=(NB2 "age" "26"); E="Women"; G="Skien";          H="Vipers Kristiansand";                                                    J="Handball";             K=(NB2 "" "Norway") },
  @{ Row=87; Id=86; B="Kristian Blummenfelt";  C="14 February 1994";   D=(NB2 "age" "26"); E="Man";                       H="Bergen Triathlon Club"; I="Arild Tveiten";                                J="Triathlon";            K=(NB2 "" "Norway") },
  @{ Row=88; Id=87; B="William Troost-Ekong";  C="1 September 1993";   D=(NB2 "age" "26"); E="Man";   G="Haarlem";        H="Udinese";                                                                 J="Football";             K=(NB2 "" "Nigeria") },
  @{ Row=89; Id=88; B="Michael Gbinije";       C="5 June 1992";        D=(NB2 "age" "27"); E="Man";   F="91 kg";          G="Hartford";                                                                J="Basketball";           K=(NB2 "" "Nigeria") },
  @{ Row=90; Id=89; B="Mahé Drysdale";         C="19 November 1978";   D="age 41";          E="Man";   F="99 kg";          G="Melbourne";          H="West End Rowing Club";                            J="Rowing";               K=(NB3 "" "New Zealand" "") },
  @{ Row=91; Id=90; B="Kelly Brazier";         C="28 October 1989";    D=(NB2 "age" "30"); E="Women"; F=(NB2 "70" "kg");  G="Dunedin";                                                                 J="Rugby";                K=(NB3 "" "New Zealand" "") },
  @{ Row=92; Id=91; B="Ferry Weertman";        C="27 June 1992";       D=(NB2 "age" "27"); E="Man";   F=(NB2 "86" "kg");  G="Naarden";                                                                 J="Swimming";             K=(NB3 "" "Netherlands" "") },
  @{ Row=93; Id=92; B="Anicka van Emden";      C="10 December 1986";   D=(NB2 "age" "33"); E="Women";                     H="Budokan Rotterdam";     I="Mark van der Ham";                             J="Judo";                 K=(NB3 "" "Netherlands" "") },
  @{ Row=94; Id=93; B="Sjef van den Berg";     C="14 April 1995";      D=(NB2 "age" "24"); E="Man";   F=(NB2 "75" "kg");  G="Heeswijk-Dinther";  H="HBV Ontspanning"; HWrap=$true; I="Ron van der Hoff"; J="Beach Volleyball";    K=(NB3 "" "Netherlands" "") }
)

foreach ($r in $rows) {
    $row = $r.Row
    $competitors.Cells.Item($row, 1).Value = $r.Id                 # A - Id
    if ($r.ContainsKey("B")) { $competitors.Cells.Item($row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $competitors.Cells.Item($row, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { $competitors.Cells.Item($row, 4).Value = $r.D }
    if ($r.ContainsKey("E")) { $competitors.Cells.Item($row, 5).Value = $r.E }
    if ($r.ContainsKey("F")) { $competitors.Cells.Item($row, 6).Value = $r.F }
    if ($r.ContainsKey("G")) { $competitors.Cells.Item($row, 7).Value = $r.G }
    if ($r.ContainsKey("H")) {
        $hCell = $competitors.Cells.Item($row, 8)
        $hCell.Value = $r.H
        if ($r.ContainsKey("HWrap") -and $r.HWrap) {
            $hCell.WrapText = $true
            $hCell.VerticalAlignment = -4108
        }
    }
    if ($r.ContainsKey("I")) { $competitors.Cells.Item($row, 9).Value = $r.I }
    if ($r.ContainsKey("J")) { $competitors.Cells.Item($row, 10).Value = $r.J }
    if ($r.ContainsKey("K")) { $competitors.Cells.Item($row, 11).Value = $r.K }
}

# Rows 95 and 96 only carry the running Id number in column A.
$competitors.Cells.Item(95, 1).Value = 94
$competitors.Cells.Item(96, 1).Value = 95

# Center-align the new "Id" cells the same way the rest of column A is
# styled (horizontal + vertical center).
$idRange = $competitors.Range("A82:A96")
$idRange.HorizontalAlignment = -4108
$idRange.VerticalAlignment = -4108

# ---------------------------------------------------------------
# 3. Minor formatting tweak on the Countries sheet: "Syria" (B90)
#    now uses left alignment.
# ---------------------------------------------------------------
$countries.Range("B90").HorizontalAlignment = -4131

# ---------------------------------------------------------------
# 4. Sheet-view bookkeeping: Competitors becomes the active sheet
#    (selection parked at D96); Countries keeps a parked selection
#    at B70 but is no longer the active tab.
# ---------------------------------------------------------------
$countries.Activate()
[void]$countries.Range("B70").Select()

$competitors.Activate()
[void]$competitors.Range("D96").Select()
